# Apply the "rich text" commit to Sheet2 ("rich test"):
#  - add a rich-text pangram in A1 (with many per-run font variations)
#  - add "hello, xssf" style combinations in column B/D, rows 3-6
#  - rename Sheet2 -> "rich test"
#  - leave the selection parked on D7 of the renamed sheet
#  - restore Sheet1 as the active/selected sheet afterwards

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------
# Row 3: "hello, xssf" in B3 (plain) and a rich "hello"(underlined) + ", xssf" in D3
# ---------------------------------------------------------------------
$ws.Range("B3").Value = "hello, xssf"

$ws.Range("D3").Value = "hello, xssf"
$ws.Range("D3").Characters(1, 5).Font.Underline = $true

# ---------------------------------------------------------------------
# Row 4: "hello, xssf" in B4 (plain) and a rich "hello, " + "xssf"(underlined) in D4
# ---------------------------------------------------------------------
$ws.Range("B4").Value = "hello, xssf"

$ws.Range("D4").Value = "hello, xssf"
$ws.Range("D4").Characters(8, 4).Font.Underline = $true

# ---------------------------------------------------------------------
# Row 5: "hello, xssf" in B5 (plain) and the SAME plain text in D5, but with the
# whole cell's font underlined (cell-level formatting, not a rich run)
# ---------------------------------------------------------------------
$ws.Range("B5").Value = "hello, xssf"

$ws.Range("D5").Value = "hello, xssf"
$ws.Range("D5").Font.Underline = $true

# ---------------------------------------------------------------------
# Row 6: "hello, xssf" in B6 (plain) and a rich "hello"(blue) + ", xssf" in D6
# ---------------------------------------------------------------------
$ws.Range("B6").Value = "hello, xssf"

$ws.Range("D6").Value = "hello, xssf"
$ws.Range("D6").Characters(1, 5).Font.Color = 12611584   # RGB(0x00,0x70,0xC0) -> 0070C0

# ---------------------------------------------------------------------
# Row 1: rich-text pangram with many different run-level font variations
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "The quick brown fox jumps over the lazy dog"

# "The" -> red
$ws.Range("A1").Characters(1, 3).Font.Color = 192                 # RGB(0xC0,0x00,0x00)

# "quick" -> bold italic
$ws.Range("A1").Characters(5, 5).Font.Bold = $true
$ws.Range("A1").Characters(5, 5).Font.Italic = $true

# " " before brown -> italic
$ws.Range("A1").Characters(10, 1).Font.Italic = $true

# "brown" -> italic, size 8
$ws.Range("A1").Characters(11, 5).Font.Italic = $true
$ws.Range("A1").Characters(11, 5).Font.Size = 8

# "jumps" -> underline, size 14, dark orange (theme Accent6, darker 50%)
$ws.Range("A1").Characters(21, 5).Font.Underline = $true
$ws.Range("A1").Characters(21, 5).Font.Size = 14
$ws.Range("A1").Characters(21, 5).Font.Color = 477336             # RGB(0x98,0x48,0x07)

# "over the lazy" -> Courier
$ws.Range("A1").Characters(27, 13).Font.Name = "Courier"

# "dog" -> bold, lighter blue (theme Text2, lighter 40%)
$ws.Range("A1").Characters(41, 3).Font.Bold = $true
$ws.Range("A1").Characters(41, 3).Font.Color = 13995605           # RGB(0x55,0x8E,0xD5)

# Row 1 is taller to fit the bigger "jumps" run
$ws.Rows.Item(1).RowHeight = 18.75

# ---------------------------------------------------------------------
# Rename the sheet and park the selection on D7, then restore Sheet1 as active
# ---------------------------------------------------------------------
$ws.Name = "rich test"

[void]$ws.Activate()
[void]$ws.Range("D7").Select()

[void]$ws1.Activate()
[void]$ws1.Range("A1").Select()
